# Excel Write + Confirmation Flow
# - Users sheet: A2 ("ID" value for Firdavs's row) becomes the literal text "Yangi Ism"
# - Sales sheet: remove the per-row Total formulas (E2:E8), including the Grand Total sum
# - Inventory sheet: remove the per-row Status formulas (E2:E6)

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Yangi Ism"

$wsSales = $wb.Worksheets.Item("Sales")
$wsSales.Range("E2:E8").ClearContents()

$wsInventory = $wb.Worksheets.Item("Inventory")
$wsInventory.Range("E2:E6").ClearContents()
